$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.488133834874986
$ws.Range("C2").Value = 0.1765239401599672
$ws.Range("D2").Value = 0.07839918899901477
$ws.Range("E2").Value = 0.08445594587336203
$ws.Range("G2").Value = 0.00245194282754721
$ws.Range("L2").Value = 0.2255077587125527
$ws.Range("N2").Value = 1.290701089977663
$ws.Range("O2").Value = 3.788337587761248

$ws.Range("B3").Value = 1.37589924506841
$ws.Range("C3").Value = 0.1629045123966932
$ws.Range("D3").Value = 0.07112518762161812
$ws.Range("E3").Value = 0.08493388216744435
$ws.Range("G3").Value = 0.002455629854554183
$ws.Range("L3").Value = 0.2161932431911424
$ws.Range("N3").Value = 1.307412781539096
$ws.Range("O3").Value = 3.766121607483427

$ws.Range("B4").Value = 1.307474604376011
$ws.Range("C4").Value = 0.1544753755962063
$ws.Range("D4").Value = 0.0666962248324694
$ws.Range("E4").Value = 0.08526774109059865
$ws.Range("G4").Value = 0.002458014952020879
$ws.Range("L4").Value = 0.2105959632905297
$ws.Range("N4").Value = 1.318227827171004
$ws.Range("O4").Value = 3.755091546260275

$ws.Range("B5").Value = 1.279714465096617
$ws.Range("C5").Value = 0.1510237399178465
$ws.Range("D5").Value = 0.06490070281617477
$ws.Range("E5").Value = 0.08541395105564398
$ws.Range("G5").Value = 0.002459017483442678
$ws.Range("L5").Value = 0.208345673034799
$ws.Range("N5").Value = 1.32277427971901
$ws.Range("O5").Value = 3.751251652438214

$ws.Range("B6").Value = 1.275112394994949
$ws.Range("C6").Value = 0.150449592890638
$ws.Range("D6").Value = 0.0646031185740128
$ws.Range("E6").Value = 0.08543884272922142
$ws.Range("G6").Value = 0.002459185803196877
$ws.Range("L6").Value = 0.2079738651207634
$ws.Range("N6").Value = 1.323537621994479
$ws.Range("O6").Value = 3.750653553803517

$ws.Range("B7").Value = 1.307099720512952
$ws.Range("C7").Value = 0.1544288931000608
$ws.Range("D7").Value = 0.06667197215655563
$ws.Range("E7").Value = 0.08526967179299305
$ws.Range("G7").Value = 0.00245802834861709
$ws.Range("L7").Value = 0.2105654909977943
$ws.Range("N7").Value = 1.318288578568747
$ws.Range("O7").Value = 3.75503711028756

$ws.Range("B8").Value = 1.449334525085931
$ws.Range("C8").Value = 0.1718418669757398
$ws.Range("D8").Value = 0.07588333017788784
$ws.Range("E8").Value = 0.084612352223008
$ws.Range("G8").Value = 0.002453189007329903
$ws.Range("L8").Value = 0.2222708167620624
$ws.Range("N8").Value = 1.296348254683618
$ws.Range("O8").Value = 3.780134743749699

$ws.Range("B9").Value = 1.732108531273127
$ws.Range("C9").Value = 0.2054565318307766
$ws.Range("D9").Value = 0.09424660527554352
$ws.Range("E9").Value = 0.08364399278820045
$ws.Range("G9").Value = 0.002444656733291866
$ws.Range("L9").Value = 0.2461933583381892
$ws.Range("N9").Value = 1.257721939652658
$ws.Range("O9").Value = 3.850142493542023

$ws.Range("B10").Value = 1.942208428356707
$ws.Range("C10").Value = 0.2298278126456808
$ws.Range("D10").Value = 0.107927588734583
$ws.Range("E10").Value = 0.08312815101600179
$ws.Range("G10").Value = 0.002438965738766132
$ws.Range("L10").Value = 0.2643638204206269
$ws.Range("N10").Value = 1.232027947307095
$ws.Range("O10").Value = 3.914371157556616

$ws.Range("B11").Value = 2.038298754259756
$ws.Range("C11").Value = 0.2408442158368018
$ws.Range("D11").Value = 0.1141940907086934
$ws.Range("E11").Value = 0.08293599352249181
$ws.Range("G11").Value = 0.002436500884019999
$ws.Range("L11").Value = 0.2727601866236569
$ws.Range("N11").Value = 1.220922831936786
$ws.Range("O11").Value = 3.94639506101106

$ws.Range("B12").Value = 2.074759279861155
$ws.Range("C12").Value = 0.2450056864387307
$ws.Range("D12").Value = 0.1165733207345596
$ws.Range("E12").Value = 0.0828693426711844
$ws.Range("G12").Value = 0.002435585240315685
$ws.Range("L12").Value = 0.2759584817303562
$ws.Range("N12").Value = 1.216801585175414
$ws.Range("O12").Value = 3.958927104605721

$ws.Range("B13").Value = 2.066903610579345
$ws.Range("C13").Value = 0.2441098947005571
$ws.Range("D13").Value = 0.1160606323398099
$ws.Range("E13").Value = 0.0828834251028745
$ws.Range("G13").Value = 0.002435781652786693
$ws.Range("L13").Value = 0.2752688366414731
$ws.Range("N13").Value = 1.217685430096257
$ws.Range("O13").Value = 3.95621005091607

$ws.Range("B14").Value = 2.041296921218304
$ws.Range("C14").Value = 0.241186787595268
$ws.Range("D14").Value = 0.1143897059295682
$ws.Range("E14").Value = 0.08293038754324122
$ws.Range("G14").Value = 0.002436425198309904
$ws.Range("L14").Value = 0.27302293601808
$ws.Range("N14").Value = 1.220582089706461
$ws.Range("O14").Value = 3.947417946935843

$ws.Range("B15").Value = 2.025621593627363
$ws.Range("C15").Value = 0.2393949683501546
$ws.Range("D15").Value = 0.1133670300435341
$ws.Range("E15").Value = 0.08295994985125965
$ws.Range("G15").Value = 0.002436821696552015
$ws.Range("L15").Value = 0.2716497020019233
$ws.Range("N15").Value = 1.222367323846207
$ws.Range("O15").Value = 3.942085362996124

$ws.Range("B16").Value = 1.935938988316877
$ws.Range("C16").Value = 0.2291064428584946
$ws.Range("D16").Value = 0.1075189279910944
$ws.Range("E16").Value = 0.08314156448400389
$ws.Range("G16").Value = 0.002439129310624848
$ws.Range("L16").Value = 0.2638177250258167
$ws.Range("N16").Value = 1.232765441472196
$ws.Range("O16").Value = 3.912334936276636

$ws.Range("B17").Value = 1.881052883170014
$ws.Range("C17").Value = 0.2227767067506647
$ws.Range("D17").Value = 0.1039423545577449
$ws.Range("E17").Value = 0.08326386719283896
$ws.Range("G17").Value = 0.002440576655031254
$ws.Range("L17").Value = 0.2590464951281888
$ws.Range("N17").Value = 1.239293832141549
$ws.Range("O17").Value = 3.894803959534727

$ws.Range("B18").Value = 1.849532354130986
$ws.Range("C18").Value = 0.2191294072664505
$ws.Range("D18").Value = 0.1018892432614251
$ws.Range("E18").Value = 0.0833382127139366
$ws.Range("G18").Value = 0.002441420806591104
$ws.Range("L18").Value = 0.2563144959283932
$ws.Range("N18").Value = 1.243103674709193
$ws.Range("O18").Value = 3.8849846079857

$ws.Range("B19").Value = 1.83886840107948
$ws.Range("C19").Value = 0.2178933635373994
$ws.Range("D19").Value = 0.1011947868008889
$ws.Range("E19").Value = 0.08336407177393035
$ws.Range("G19").Value = 0.002441708630041346
$ws.Range("L19").Value = 0.2553915979980701
$ws.Range("N19").Value = 1.24440304387527
$ws.Range("O19").Value = 3.881705232542799

$ws.Range("B20").Value = 1.886890588495248
$ws.Range("C20").Value = 0.2234512020868635
$ws.Range("D20").Value = 0.1043226684004139
$ws.Range("E20").Value = 0.08325043382718533
$ws.Range("G20").Value = 0.002440421374558749
$ws.Range("L20").Value = 0.2595531291185154
$ws.Range("N20").Value = 1.238593192319435
$ws.Range("O20").Value = 3.896642825689071

$ws.Range("B21").Value = 2.048816246389265
$ws.Range("C21").Value = 0.2420456527297858
$ws.Range("D21").Value = 0.1148803276025632
$ws.Range("E21").Value = 0.08291642754887896
$ws.Range("G21").Value = 0.002436235692771067
$ws.Range("L21").Value = 0.273682101937311
$ws.Range("N21").Value = 1.219728989040465
$ws.Range("O21").Value = 3.949989387306971

$ws.Range("B22").Value = 2.155070340114833
$ws.Range("C22").Value = 0.2541386866430173
$ws.Range("D22").Value = 0.1218167800271459
$ws.Range("E22").Value = 0.0827337797872616
$ws.Range("G22").Value = 0.002433603488303742
$ws.Range("L22").Value = 0.2830256305401377
$ws.Range("N22").Value = 1.207889946892493
$ws.Range("O22").Value = 3.987217607361458

$ws.Range("B23").Value = 2.098321793416517
$ws.Range("C23").Value = 0.2476898870200159
$ws.Range("D23").Value = 0.1181113116305568
$ws.Range("E23").Value = 0.08282799971327393
$ws.Range("G23").Value = 0.002434998915012945
$ws.Range("L23").Value = 0.2780287995440744
$ws.Range("N23").Value = 1.214163796782895
$ws.Range("O23").Value = 3.967131381702188

$ws.Range("B24").Value = 1.884251255250661
$ws.Range("C24").Value = 0.2231462884388691
$ws.Range("D24").Value = 0.1041507188379427
$ws.Range("E24").Value = 0.08325649449277783
$ws.Range("G24").Value = 0.0024404915392942
$ws.Range("L24").Value = 0.2593240454923631
$ws.Range("N24").Value = 1.238909775328743
$ws.Range("O24").Value = 3.895810666268062

$ws.Range("B25").Value = 1.655198578221757
$ws.Range("C25").Value = 0.1964199503549366
$ws.Range("D25").Value = 0.08924606170805305
$ws.Range("E25").Value = 0.08387161863626069
$ws.Range("G25").Value = 0.002446863053302665
$ws.Range("L25").Value = 0.239617569831637
$ws.Range("N25").Value = 1.267700166927771
$ws.Range("O25").Value = 3.828965420608853
